$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.13508
$ws.Range("H2").Value = 3.40524
$ws.Range("I2").Value = 0.1224366388308639
$ws.Range("J2").Value = 0.1224366388308639
$ws.Range("Q2").Value = 0.04569037524000001
$ws.Range("R2").Value = 0.41121337716
$ws.Range("S2").Value = 0.1224366388308639
$ws.Range("T2").Value = 0.1224366388308639

# Row 3 updates
$ws.Range("G3").Value = 0.672624
$ws.Range("H3").Value = 2.017872
$ws.Range("I3").Value = 0.07255331937570129
$ws.Range("J3").Value = 0.07255331937570129
$ws.Range("Q3").Value = 0.027075133872
$ws.Range("R3").Value = 0.243676204848
$ws.Range("S3").Value = 0.07255331937570129
$ws.Range("T3").Value = 0.07255331937570129

# Row 4 updates
$ws.Range("G4").Value = 7.463050333333332
$ws.Range("H4").Value = 22.389151
$ws.Range("I4").Value = 0.8050100417934347
$ws.Range("J4").Value = 0.8050100417934348
$ws.Range("Q4").Value = 0.3004101650676667
$ws.Range("R4").Value = 2.703691485609
$ws.Range("S4").Value = 0.8050100417934347
$ws.Range("T4").Value = 0.8050100417934348
